$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update utilisation values in column E (rows 2-6): 600000 -> 800000
$ws.Range("E2").Value = 800000
$ws.Range("E3").Value = 800000
$ws.Range("E4").Value = 800000
$ws.Range("E5").Value = 800000
$ws.Range("E6").Value = 800000

# Swap G5 / H5 values: G5 was 3938753.8, H5 was 456 -> G5 becomes 456, H5 becomes 3938753.8
$ws.Range("G5").Value = 456
$ws.Range("H5").Value = 3938753.8

# Update the active selection to match the new state (E2:H6, active cell E2)
$ws.Range("E2:H6").Select()
